$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: Insert 54 blank rows before the old row 1185 (Day 6 header), pushing it and everything after down by 54
$ws.Rows("1185:1238").Insert()

# Step 2: Insert 318 more blank rows right after the relocated Day 6 header (now at row 1239),
# pushing the Day 7 section (now at 1244) down to 1562
$ws.Rows("1240:1557").Insert()

# Step 3: Fill in the new Day 6 dialogue content
$ws.Range("A1239").Value = '//'
$ws.Range("B1239").Value = 'Day 6'
$ws.Range("A1240").Value = 'CONVERSATION'
$ws.Range("B1240").Value = 'Check TV'
$ws.Range("D1240").Value = 6
$ws.Range("A1241").Value = 'TV'
$ws.Range("B1241").Value = '*STATIC*'
$ws.Range("D1241").Value = 6
$ws.Range("E1241").Value = 'tv6'
$ws.Range("F1241").Value = 'DRAFT'
$ws.Range("G1241").Value = 'NONE'
$ws.Range("A1242").Value = 'END'
$ws.Range("A1245").Value = 'CONVERSATION'
$ws.Range("B1245").Value = 'Check Radio'
$ws.Range("D1245").Value = 6
$ws.Range("A1246").Value = 'Radio'
$ws.Range("B1246").Value = 'The war in Brasnia continues to expand outwards. Citizens of the north and east are expected to hear firefights and artillery. Ceasefire negotiations continue with the US still apprehensive.'
$ws.Range("D1246").Value = 6
$ws.Range("E1246").Value = 'r6'
$ws.Range("F1246").Value = 'DRAFT'
$ws.Range("G1246").Value = 'NONE'
$ws.Range("A1247").Value = 'END'
$ws.Range("A1249").Value = 'CONVERSATION'
$ws.Range("B1249").Value = 'Check Outside'
$ws.Range("D1249").Value = 6
$ws.Range("A1250").Value = 'Player'
$ws.Range("B1250").Value = 'It looks clearer and clearer by the day... still hostile though.'
$ws.Range("D1250").Value = 6
$ws.Range("E1250").Value = 'p6_outside'
$ws.Range("F1250").Value = 'DRAFT'
$ws.Range("G1250").Value = 'NONE'
$ws.Range("A1251").Value = 'END'
$ws.Range("A1253").Value = 'CONVERSATION'
$ws.Range("B1253").Value = 'Check Supplies'
$ws.Range("D1253").Value = 6
$ws.Range("A1254").Value = 'Player'
$ws.Range("B1254").Value = 'We should be alright for a few days at least.'
$ws.Range("D1254").Value = 6
$ws.Range("E1254").Value = 'p6_supplies'
$ws.Range("F1254").Value = 'DRAFT'
$ws.Range("G1254").Value = 'NONE'
$ws.Range("A1255").Value = 'END'
$ws.Range("A1260").Value = 'CONVERSATION'
$ws.Range("B1260").Value = 'Bob wakes up'
$ws.Range("D1260").Value = 6
$ws.Range("A1261").Value = 'Bob'
$ws.Range("B1261").Value = 'Big guys got the TV down by the looks of it.'
$ws.Range("D1261").Value = 6
$ws.Range("E1261").Value = 'b_6_tv_taken_down'
$ws.Range("F1261").Value = 'DRAFT'
$ws.Range("G1261").Value = 'NONE'
$ws.Range("A1262").Value = 'Player'
$ws.Range("B1262").Value = 'We could use some more information on what''s happening.'
$ws.Range("D1262").Value = 6
$ws.Range("E1262").Value = 'p_6_tv_could_use_info'
$ws.Range("F1262").Value = 'DRAFT'
$ws.Range("G1262").Value = 'NONE'
$ws.Range("A1263").Value = 'Bob'
$ws.Range("B1263").Value = 'It''s just fighting and more fighting.'
$ws.Range("D1263").Value = 6
$ws.Range("E1263").Value = 'b_6_tv_fighting'
$ws.Range("F1263").Value = 'DRAFT'
$ws.Range("G1263").Value = 'NONE'
$ws.Range("A1264").Value = 'END'
$ws.Range("A1266").Value = 'CONVERSATION'
$ws.Range("B1266").Value = 'Jessica wakes up'
$ws.Range("D1266").Value = 6
$ws.Range("A1267").Value = 'Jessica'
$ws.Range("B1267").Value = 'Good morning everyone!'
$ws.Range("D1267").Value = 6
$ws.Range("E1267").Value = 'j_6_good_morning'
$ws.Range("F1267").Value = 'DRAFT'
$ws.Range("G1267").Value = 'NONE'
$ws.Range("A1268").Value = 'IF'
$ws.Range("B1268").Value = 'Bob'
$ws.Range("A1269").Value = 'Bob'
$ws.Range("B1269").Value = '*sigh*'
$ws.Range("D1269").Value = 6
$ws.Range("E1269").Value = 'b_6_sigh'
$ws.Range("F1269").Value = 'DRAFT'
$ws.Range("G1269").Value = 'NONE'
$ws.Range("A1270").Value = 'ELSE'
$ws.Range("A1271").Value = 'Player'
$ws.Range("B1271").Value = 'Hi.'
$ws.Range("D1271").Value = 6
$ws.Range("E1271").Value = 'p_6_j_hi'
$ws.Range("F1271").Value = 'DRAFT'
$ws.Range("G1271").Value = 'NONE'
$ws.Range("A1272").Value = 'ENDIF'
$ws.Range("A1273").Value = 'END'
$ws.Range("A1275").Value = 'CONVERSATION'
$ws.Range("B1275").Value = 'Violet wakes up'
$ws.Range("A1276").Value = 'Violet'
$ws.Range("B1276").Value = '*your beds arent half bad actually*'
$ws.Range("D1276").Value = 6
$ws.Range("E1276").Value = 'v_6_wake_up'
$ws.Range("F1276").Value = 'DRAFT'
$ws.Range("G1276").Value = 'NONE'
$ws.Range("A1277").Value = 'END'
$ws.Range("A1279").Value = 'CONVERSATION'
$ws.Range("B1279").Value = 'Dad wakes up'
$ws.Range("D1279").Value = 6
$ws.Range("A1280").Value = 'Dad'
$ws.Range("B1280").Value = 'Morning guys... Thanks for letting me stay the night.'
$ws.Range("D1280").Value = 6
$ws.Range("E1280").Value = 'dad_6_wake_up'
$ws.Range("F1280").Value = 'DRAFT'
$ws.Range("G1280").Value = 'NONE'
$ws.Range("A1281").Value = 'IF'
$ws.Range("B1281").Value = 'Bob'
$ws.Range("A1282").Value = 'Bob'
$ws.Range("B1282").Value = 'Don''t want to be that guy, but you did say you were just staying the night.'
$ws.Range("D1282").Value = 6
$ws.Range("E1282").Value = 'b_6_d_you_gotta_leave'
$ws.Range("F1282").Value = 'DRAFT'
$ws.Range("G1282").Value = 'NONE'
$ws.Range("A1283").Value = 'ENDIF'
$ws.Range("A1284").Value = 'Dad'
$ws.Range("B1284").Value = 'I''m a man of my word. It was nice getting to briefly meet you folks, but it''s time for me to hit the road... My son''s out there somewhere.'
$ws.Range("D1284").Value = 6
$ws.Range("E1284").Value = 'dad_6_i_will_leave'
$ws.Range("F1284").Value = 'DRAFT'
$ws.Range("G1284").Value = 'NONE'
$ws.Range("A1285").Value = 'IF'
$ws.Range("B1285").Value = 'Jessica'
$ws.Range("A1286").Value = 'Jessica'
$ws.Range("B1286").Value = 'Is there anything we can give you to help you out?'
$ws.Range("D1286").Value = 6
$ws.Range("E1286").Value = 'j_6_d_anything_to_give'
$ws.Range("F1286").Value = 'DRAFT'
$ws.Range("G1286").Value = 'NONE'
$ws.Range("A1287").Value = 'Dad'
$ws.Range("B1287").Value = 'I''m not sure... I''ve got food, water, clothes on my back and a good pair of mitts.'
$ws.Range("D1287").Value = 6
$ws.Range("E1287").Value = 'dad_6_dont_need_anything'
$ws.Range("F1287").Value = 'DRAFT'
$ws.Range("G1287").Value = 'NONE'
$ws.Range("A1288").Value = 'ENDIF'
$ws.Range("A1289").Value = 'IF'
$ws.Range("B1289").Value = '!HasShotgun'
$ws.Range("A1290").Value = 'GOTO'
$ws.Range("B1290").Value = 'Nothing for Dad'
$ws.Range("A1291").Value = 'ENDIF'
$ws.Range("A1292").Value = 'CHOICE'
$ws.Range("B1292").Value = 'Offer Dad shotgun'
$ws.Range("A1293").Value = 'Player'
$ws.Range("B1293").Value = 'Offer Shotgun'
$ws.Range("D1293").Value = 6
$ws.Range("E1293").Value = 'p_6_d_offer_shotgun'
$ws.Range("F1293").Value = 'DRAFT'
$ws.Range("G1293").Value = 'NONE'
$ws.Range("A1294").Value = 'CHOICE'
$ws.Range("B1294").Value = 'Nothing for Dad'
$ws.Range("A1295").Value = 'Player'
$ws.Range("B1295").Value = 'Send him on his way'
$ws.Range("D1295").Value = 6
$ws.Range("E1295").Value = 'p_6_d_send_dad_along'
$ws.Range("F1295").Value = 'DRAFT'
$ws.Range("G1295").Value = 'NONE'
$ws.Range("A1296").Value = 'END'
$ws.Range("A1299").Value = 'CONVERSATION'
$ws.Range("B1299").Value = 'Nothing for Dad'
$ws.Range("D1299").Value = 6
$ws.Range("A1300").Value = 'Player'
$ws.Range("B1300").Value = 'Sorry, but I don''t think there''s anything else we can do for you.'
$ws.Range("D1300").Value = 6
$ws.Range("E1300").Value = 'p_6_d_nothing_else'
$ws.Range("F1300").Value = 'DRAFT'
$ws.Range("G1300").Value = 'NONE'
$ws.Range("A1301").Value = 'Dad'
$ws.Range("B1301").Value = 'That''s alright, just spending the night is nice enough.'
$ws.Range("D1301").Value = 6
$ws.Range("E1301").Value = 'dad_6_nothing_is_alright'
$ws.Range("F1301").Value = 'DRAFT'
$ws.Range("G1301").Value = 'NONE'
$ws.Range("A1302").Value = 'GOTO'
$ws.Range("B1302").Value = 'Dad leaves'
$ws.Range("A1303").Value = 'END'
$ws.Range("A1305").Value = 'CONVERSATION'
$ws.Range("B1305").Value = 'Offer Dad shotgun'
$ws.Range("D1305").Value = 6
$ws.Range("A1306").Value = 'Dad'
$ws.Range("B1306").Value = 'It would certainly help me out.'
$ws.Range("D1306").Value = 6
$ws.Range("E1306").Value = 'dad_6_shotgun_would_help'
$ws.Range("F1306").Value = 'DRAFT'
$ws.Range("G1306").Value = 'NONE'
$ws.Range("A1307").Value = 'IF'
$ws.Range("B1307").Value = 'Bob && !Jessica'
$ws.Range("A1308").Value = 'Bob'
$ws.Range("B1308").Value = '(to you) Don''t give him our gun. He''s only been here for a night and he could still be dangerous.'
$ws.Range("D1308").Value = 6
$ws.Range("E1308").Value = 'b_6_d_dont_give_gun'
$ws.Range("F1308").Value = 'DRAFT'
$ws.Range("G1308").Value = 'NONE'
$ws.Range("A1309").Value = 'ELIF'
$ws.Range("B1309").Value = 'Bob && Jessica'
$ws.Range("A1310").Value = 'Bob'
$ws.Range("B1310").Value = '(to you) Be careful with him...'
$ws.Range("D1310").Value = 6
$ws.Range("E1310").Value = 'b_6_d_be_cautious'
$ws.Range("F1310").Value = 'DRAFT'
$ws.Range("G1310").Value = 'NONE'
$ws.Range("A1311").Value = 'ENDIF'
$ws.Range("A1312").Value = 'IF'
$ws.Range("B1312").Value = 'Jessica'
$ws.Range("A1313").Value = 'Jessica'
$ws.Range("B1313").Value = '(to you) Give it to him, he can use it a lot more than we can.'
$ws.Range("D1313").Value = 6
$ws.Range("E1313").Value = 'j_6_d_give_gun'
$ws.Range("F1313").Value = 'DRAFT'
$ws.Range("G1313").Value = 'NONE'
$ws.Range("A1314").Value = 'ENDIF'
$ws.Range("A1315").Value = 'IF'
$ws.Range("B1315").Value = 'Violet'
$ws.Range("A1316").Value = 'Violet'
$ws.Range("B1316").Value = '*ask him some more questions before you give him the gun*'
$ws.Range("D1316").Value = 6
$ws.Range("E1316").Value = 'v_6_d_ask_questions_before'
$ws.Range("F1316").Value = 'DRAFT'
$ws.Range("G1316").Value = 'NONE'
$ws.Range("A1317").Value = 'ENDIF'
$ws.Range("A1318").Value = 'GOTO'
$ws.Range("B1318").Value = 'Dad questions'
$ws.Range("A1319").Value = 'END'
$ws.Range("A1322").Value = 'CONVERSATION'
$ws.Range("B1322").Value = 'Dad questions'
$ws.Range("D1322").Value = 6
$ws.Range("A1323").Value = 'CHOICE'
$ws.Range("B1323").Value = 'Will Dad bring shotgun back'
$ws.Range("A1324").Value = 'Player'
$ws.Range("B1324").Value = '"Will you bring the gun back to us?"'
$ws.Range("D1324").Value = 6
$ws.Range("E1324").Value = 'p_6_d_will_you_bring_back'
$ws.Range("F1324").Value = 'DRAFT'
$ws.Range("G1324").Value = 'NONE'
$ws.Range("A1325").Value = 'CHOICE'
$ws.Range("B1325").Value = 'Will shotgun help Dad'
$ws.Range("A1326").Value = 'Player'
$ws.Range("B1326").Value = '"Do you really think this would help you out?"'
$ws.Range("D1326").Value = 6
$ws.Range("E1326").Value = 'p_6_d_will_gun_help'
$ws.Range("F1326").Value = 'DRAFT'
$ws.Range("G1326").Value = 'NONE'
$ws.Range("A1327").Value = 'CHOICE'
$ws.Range("B1327").Value = 'Give Dad shotgun'
$ws.Range("A1328").Value = 'Player'
$ws.Range("B1328").Value = 'Give Shotgun'
$ws.Range("D1328").Value = 6
$ws.Range("E1328").Value = 'p_6_d_give_shotgun'
$ws.Range("F1328").Value = 'DRAFT'
$ws.Range("G1328").Value = 'NONE'
$ws.Range("A1329").Value = 'CHOICE'
$ws.Range("B1329").Value = 'Keep shotgun'
$ws.Range("A1330").Value = 'Player'
$ws.Range("B1330").Value = 'Keep Shotgun'
$ws.Range("D1330").Value = 6
$ws.Range("E1330").Value = 'p_6_d_keep_shotgun'
$ws.Range("F1330").Value = 'DRAFT'
$ws.Range("G1330").Value = 'NONE'
$ws.Range("A1331").Value = 'END'
$ws.Range("A1333").Value = 'CONVERSATION'
$ws.Range("B1333").Value = 'Will Dad bring shotgun back'
$ws.Range("D1333").Value = 6
$ws.Range("A1334").Value = 'Dad'
$ws.Range("B1334").Value = 'I''m not sure if you would be seeing it again. I don''t plan on coming back, but who knows what way the wind blows.'
$ws.Range("D1334").Value = 6
$ws.Range("E1334").Value = 'dad_6_p_wont_bring_back'
$ws.Range("F1334").Value = 'DRAFT'
$ws.Range("G1334").Value = 'NONE'
$ws.Range("A1335").Value = 'GOTO'
$ws.Range("B1335").Value = 'Dad questions'
$ws.Range("A1336").Value = 'END'
$ws.Range("A1338").Value = 'CONVERSATION'
$ws.Range("B1338").Value = 'Will shotgun help Dad'
$ws.Range("D1338").Value = 6
$ws.Range("A1339").Value = 'Dad'
$ws.Range("B1339").Value = 'Yes, my son might be somewhere dangerous.'
$ws.Range("D1339").Value = 6
$ws.Range("E1339").Value = 'dad_6_p_son_maybe_danger'
$ws.Range("F1339").Value = 'DRAFT'
$ws.Range("G1339").Value = 'NONE'
$ws.Range("A1340").Value = 'GOTO'
$ws.Range("B1340").Value = 'Dad questions'
$ws.Range("A1341").Value = 'END'
$ws.Range("A1343").Value = 'CONVERSATION'
$ws.Range("B1343").Value = 'Give Dad shotgun'
$ws.Range("D1343").Value = 6
$ws.Range("A1344").Value = 'IF'
$ws.Range("B1344").Value = 'Bob'
$ws.Range("A1345").Value = 'Bob'
$ws.Range("B1345").Value = '*sigh*'
$ws.Range("D1345").Value = 6
$ws.Range("E1345").Value = 'b_6_sigh'
$ws.Range("F1345").Value = 'DRAFT'
$ws.Range("G1345").Value = 'NONE'
$ws.Range("A1346").Value = 'ENDIF'
$ws.Range("A1347").Value = 'Dad'
$ws.Range("B1347").Value = 'Thank you... I will remember this.'
$ws.Range("D1347").Value = 6
$ws.Range("E1347").Value = 'dad_6_thanks_for_gun'
$ws.Range("F1347").Value = 'DRAFT'
$ws.Range("G1347").Value = 'NONE'
$ws.Range("A1348").Value = 'GOTO'
$ws.Range("B1348").Value = 'Dad leaves'
$ws.Range("A1349").Value = 'END'
$ws.Range("A1351").Value = 'CONVERSATION'
$ws.Range("B1351").Value = 'Keep shotgun'
$ws.Range("D1351").Value = 6
$ws.Range("A1352").Value = 'Player'
$ws.Range("B1352").Value = 'I think we should keep the shotgun.'
$ws.Range("D1352").Value = 6
$ws.Range("E1352").Value = 'p_6_d_we_should_keep_gun'
$ws.Range("F1352").Value = 'DRAFT'
$ws.Range("G1352").Value = 'NONE'
$ws.Range("A1353").Value = 'Dad'
$ws.Range("B1353").Value = 'I understand. You have to keep yourself safe too.'
$ws.Range("D1353").Value = 6
$ws.Range("E1353").Value = 'dad_6_understand'
$ws.Range("F1353").Value = 'DRAFT'
$ws.Range("G1353").Value = 'NONE'
$ws.Range("A1354").Value = 'GOTO'
$ws.Range("B1354").Value = 'Dad leaves'
$ws.Range("A1355").Value = 'END'
$ws.Range("A1361").Value = 'CONVERSATION'
$ws.Range("B1361").Value = 'Dad leaves'
$ws.Range("D1361").Value = 6
$ws.Range("A1362").Value = 'Dad'
$ws.Range("B1362").Value = 'Regardless, thanks for the hospitality. Stay safe!'
$ws.Range("D1362").Value = 6
$ws.Range("E1362").Value = 'dad_6_goodbye'
$ws.Range("A1363").Value = 'IF'
$ws.Range("B1363").Value = 'Jessica'
$ws.Range("A1364").Value = 'Jessica'
$ws.Range("B1364").Value = 'Bye! Hope you can find your son!'
$ws.Range("D1364").Value = 6
$ws.Range("E1364").Value = 'j_6_d_goodbye'
$ws.Range("A1365").Value = 'ENDIF'
$ws.Range("A1366").Value = 'IF'
$ws.Range("B1366").Value = 'Bob'
$ws.Range("A1367").Value = 'Bob'
$ws.Range("B1367").Value = 'Good luck.'
$ws.Range("D1367").Value = 6
$ws.Range("E1367").Value = 'b_6_d_goodbye'
$ws.Range("A1368").Value = 'ENDIF'
$ws.Range("A1369").Value = 'IF'
$ws.Range("B1369").Value = 'Violet'
$ws.Range("A1370").Value = 'Violet'
$ws.Range("B1370").Value = '(to you) *i don''t think he''ll find his son*'
$ws.Range("D1370").Value = 6
$ws.Range("E1370").Value = 'v_6_d_goodbye'
$ws.Range("A1371").Value = 'ENDIF'
$ws.Range("A1372").Value = 'END'
